$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "24.659.33"
$ws.Range("E2").Value = "  -0.25%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.676.73"
$ws.Range("E3").Value = "  -0.38%  "

$ws.Range("E4").Value = "  +0.38%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "307.15"
$ws.Range("E5").Value = "  +0.34%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9993"
$ws.Range("E6").Value = "  +0.39%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3686"
$ws.Range("E7").Value = "  +0.11%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "48.26"
$ws.Range("E8").Value = "  -1.50%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3371"
$ws.Range("E9").Value = "  -1.39%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.177"
$ws.Range("E10").Value = "  +1.44%  "

$ws.Range("E11").Value = "  +1.38%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.000"
$ws.Range("E12").Value = "  +0.31%  "

$ws.Range("B13").Value = "Solana"
$ws.Range("C13").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.54"
$ws.Range("E13").Value = "  +2.13%  "

$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.183"
$ws.Range("E14").Value = "  +1.64%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.799"
$ws.Range("E15").Value = "  +1.72%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.673.97"
$ws.Range("E16").Value = "  -0.36%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001098"
$ws.Range("E17").Value = "  -0.31%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06636"
$ws.Range("E18").Value = "  -0.16%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.9992"
$ws.Range("E19").Value = "  +0.44%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "81.67"
$ws.Range("E20").Value = "  +1.21%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "16.82"
$ws.Range("E21").Value = "  +2.85%  "

$ws.Range("E22").Value = "  +2.19%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.69"
$ws.Range("E23").Value = "  +5.08%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "24.621.72"
$ws.Range("E24").Value = "  -0.10%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.442"
$ws.Range("E25").Value = "  +1.60%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.696"
$ws.Range("E26").Value = "  +1.57%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.85"
$ws.Range("E27").Value = "  +2.19%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "149.19"
$ws.Range("E28").Value = "  -2.21%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "130.07"
$ws.Range("E29").Value = "  +2.01%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.857.80"
$ws.Range("E30").Value = "  -0.46%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.216"
$ws.Range("E31").Value = "  +24.85%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.515"
$ws.Range("E32").Value = "  +4.40%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.155"
$ws.Range("E33").Value = "  +3.34%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08591"
$ws.Range("E34").Value = "  +2.27%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "13.31"
$ws.Range("E35").Value = "  +8.30%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.715"
$ws.Range("E36").Value = "  +1.93%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.420"
$ws.Range("E37").Value = "  +2.34%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.06461"
$ws.Range("E38").Value = "  +2.01%  "

$ws.Range("B39").Value = "FraxShare"
$ws.Range("C39").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.848"
$ws.Range("E39").Value = "  +2.78%  "

$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.02354"
$ws.Range("E40").Value = "  +2.23%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.2167"
$ws.Range("E41").Value = "  +4.14%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.235"
$ws.Range("E42").Value = "  -0.60%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.6259"
$ws.Range("E43").Value = "  +3.02%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.0000"
$ws.Range("E44").Value = "  +0.50%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.48"
$ws.Range("E45").Value = "  +3.68%  "

$ws.Range("E46").Value = "  +0.55%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5943"
$ws.Range("E47").Value = "  +1.25%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.045"
$ws.Range("E48").Value = "  +2.51%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "126.19"
$ws.Range("E49").Value = "  +0.67%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.07132"
$ws.Range("E50").Value = "  -1.22%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "77.09"
$ws.Range("E51").Value = "  +1.98%  "
